{"js": "// Replace the 100 arithmetic-problem cells (20 rows x 5 cols) in the\n// single table of the document with their new values, in document order.\n// Built directly from the authoritative old->new diff; values are written\n// positionally (row-major) so duplicate problem text (e.g. \"52-6=\", which\n// appears twice in the source with two different replacements) is handled\n// correctly.\nconst newValues = [\n  [\"81-69=\", \"24+39=\", \"90-18=\", \"71-62=\", \"53-46=\"],\n  [\"8+13=\", \"77-38=\", \"14+83=\", \"0+54=\", \"64+17=\"],\n  [\"53-17=\", \"36-13=\", \"65+15=\", \"9-0=\", \"57+17=\"],\n  [\"63-63=\", \"42-18=\", \"85-74=\", \"2+89=\", \"70-39=\"],\n  [\"60+34=\", \"67-46=\", \"90-55=\", \"34-33=\", \"78-35=\"],\n  [\"41-17=\", \"44-15=\", \"86-26=\", \"56+22=\", \"90-78=\"],\n  [\"52-45=\", \"40+5=\", \"78-12=\", \"43-20=\", \"77-43=\"],\n  [\"16+31=\", \"12+59=\", \"13+33=\", \"44-37=\", \"87-30=\"],\n  [\"14+56=\", \"19+48=\", \"78-54=\", \"72-0=\", \"55-12=\"],\n  [\"77-76=\", \"18+33=\", \"93-72=\", \"49+13=\", \"43+42=\"],\n  [\"30-9=\", \"44-1=\", \"69-27=\", \"52-14=\", \"99-67=\"],\n  [\"98-63=\", \"34+56=\", \"1+89=\", \"97-29=\", \"49-32=\"],\n  [\"58-57=\", \"21+16=\", \"38+19=\", \"90-39=\", \"45+18=\"],\n  [\"60-59=\", \"96-28=\", \"0+59=\", \"8+49=\", \"58-38=\"],\n  [\"2+1=\", \"66+33=\", \"21+5=\", \"56-34=\", \"65-35=\"],\n  [\"30+44=\", \"27-16=\", \"54-47=\", \"97-23=\", \"8+88=\"],\n  [\"13+23=\", \"82-58=\", \"17+12=\", \"83-78=\", \"6+81=\"],\n  [\"17+43=\", \"54-40=\", \"48-18=\", \"70-10=\", \"28+0=\"],\n  [\"36+12=\", \"78-37=\", \"5+28=\", \"27+3=\", \"49-44=\"],\n  [\"94-87=\", \"25+47=\", \"28+47=\", \"68+25=\", \"24+6=\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.rowCount !== newValues.length) {\n  throw new Error(`Expected ${newValues.length} rows, found ${table.rowCount}.`);\n}\n\n// Table.values setter (Office.js) walks the 2D array and assigns\n// getCell(r, c).value = v[r][c] for each cell, which replaces the cell's\n// text while preserving the existing run/paragraph formatting (font,\n// size, alignment, etc.).\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-problem cells (20 rows x 5 cols) in the\n# single table of the document with their new values, in document order.\n# Built directly from the authoritative old->new diff; values are written\n# positionally (row-major) via Cell(row, col), which is robust to\n# duplicate problem text (e.g. \"52-6=\" appears twice in the source).\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$newValues = @(\n    @('81-69=', '24+39=', '90-18=', '71-62=', '53-46='),\n    @('8+13=', '77-38=', '14+83=', '0+54=', '64+17='),\n    @('53-17=', '36-13=', '65+15=', '9-0=', '57+17='),\n    @('63-63=', '42-18=', '85-74=', '2+89=', '70-39='),\n    @('60+34=', '67-46=', '90-55=', '34-33=', '78-35='),\n    @('41-17=', '44-15=', '86-26=', '56+22=', '90-78='),\n    @('52-45=', '40+5=', '78-12=', '43-20=', '77-43='),\n    @('16+31=', '12+59=', '13+33=', '44-37=', '87-30='),\n    @('14+56=', '19+48=', '78-54=', '72-0=', '55-12='),\n    @('77-76=', '18+33=', '93-72=', '49+13=', '43+42='),\n    @('30-9=', '44-1=', '69-27=', '52-14=', '99-67='),\n    @('98-63=', '34+56=', '1+89=', '97-29=', '49-32='),\n    @('58-57=', '21+16=', '38+19=', '90-39=', '45+18='),\n    @('60-59=', '96-28=', '0+59=', '8+49=', '58-38='),\n    @('2+1=', '66+33=', '21+5=', '56-34=', '65-35='),\n    @('30+44=', '27-16=', '54-47=', '97-23=', '8+88='),\n    @('13+23=', '82-58=', '17+12=', '83-78=', '6+81='),\n    @('17+43=', '54-40=', '48-18=', '70-10=', '28+0='),\n    @('36+12=', '78-37=', '5+28=', '27+3=', '49-44='),\n    @('94-87=', '25+47=', '28+47=', '68+25=', '24+6=')\n)\n\nif ($table.Rows.Count -ne $newValues.Count) {\n    throw \"Expected $($newValues.Count) rows, found $($table.Rows.Count).\"\n}\n\nfor ($r = 0; $r -lt $newValues.Count; $r++) {\n    $rowValues = $newValues[$r]\n    for ($c = 0; $c -lt $rowValues.Count; $c++) {\n        # Word COM is 1-based: Cell(row, col).\n        $table.Cell($r + 1, $c + 1).Range.Text = $rowValues[$c]\n    }\n}\n"}
